$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Date value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2026-01-23T08:28:04+00:00"

# --- Mapping Table 1 sheet: update relationship mapping row ---
$map1 = $wb.Worksheets.Item("Mapping Table 1")
$map1.Range("A5").Value = "FRCDAAntecedentsFamiliaux.subject"
$map1.Range("D5").Value = "FRFamilyMemberHistoryDocument.relationship"
